$d = $word.ActiveDocument
$startPar = $d.Paragraphs.Item(61)
$endPar = $d.Paragraphs.Item(66)
$r = $d.Range($startPar.Range.Start, $endPar.Range.End)
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>{{OWNERSHIP}}</w:t></w:r>
<w:bookmarkStart w:id="100" w:name="_6o8u7emblwbs" w:colFirst="0" w:colLast="0"/>
<w:bookmarkStart w:id="101" w:name="_GoBack"/>
<w:bookmarkEnd w:id="100"/>
<w:bookmarkEnd w:id="101"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$r.InsertXML($xml)
